$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new columns ("ASCTB_unique_CT" and "Azimuth_unique_CT") right
# after "Az_Asctb_perfect_matches" (column D), pushing the former E:J columns
# over to G:L.
# ---------------------------------------------------------------------------
$ws.Columns("E:F").Insert()

# New column headers
$ws.Range("E1").Value = "ASCTB_unique_CT"
$ws.Range("F1").Value = "Azimuth_unique_CT"

# The two "match_found_corsswalk" headers (now shifted to I1/J1) were renamed
# to "*_ct_match_found_corsswalk"
$ws.Range("I1").Value = "Az_ct_match_found_corsswalk"
$ws.Range("J1").Value = "Asctb_ct_match_found_corsswalk"

# ---------------------------------------------------------------------------
# Row data (organ rows 2-7). Columns B:D and G:L keep their old numbers
# (shifted right by two columns); columns E:F hold the new unique-CT counts.
# ---------------------------------------------------------------------------

# Row 2 - lung
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 14
$ws.Range("E2").Value = 64
$ws.Range("F2").Value = 43
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 31
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 28
$ws.Range("K2").Value = 22
$ws.Range("L2").Value = 25

# Row 3 - pancreas
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 9
$ws.Range("E3").Value = 26
$ws.Range("F3").Value = 12
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = 12
$ws.Range("L3").Value = 3

# Row 4 - kidney
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 12
$ws.Range("D4").Value = 26
$ws.Range("E4").Value = 47
$ws.Range("F4").Value = 42
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 12
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 11
$ws.Range("K4").Value = 10
$ws.Range("L4").Value = 14

# Row 5 - brain
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 10

# Row 6 - bone_marrow
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 20
$ws.Range("E6").Value = 37
$ws.Range("F6").Value = 38
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 9
$ws.Range("J6").Value = 6
$ws.Range("K6").Value = 11
$ws.Range("L6").Value = 9

# Row 7 - blood_pmbc
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 11
$ws.Range("E7").Value = 26
$ws.Range("F7").Value = 36
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 16
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 13
$ws.Range("L7").Value = 9

# ---------------------------------------------------------------------------
# Cosmetic touch-ups matching the author's final save (column widths were
# re-autofit after the insert, and the active selection moved to D5).
# ---------------------------------------------------------------------------
$ws.Columns(1).ColumnWidth = 11.830729166666666
$ws.Columns(2).ColumnWidth = 13.053385416666666
$ws.Columns(3).ColumnWidth = 16.385416666666668
$ws.Columns(4).ColumnWidth = 24.053385416666668
$ws.Columns(5).ColumnWidth = 15.608072916666666
$ws.Columns(6).ColumnWidth = 17.944010416666668
$ws.Columns(7).ColumnWidth = 20.276041666666668
$ws.Columns(8).ColumnWidth = 21.608072916666668
$ws.Columns(9).ColumnWidth = 26.830729166666668
$ws.Columns(10).ColumnWidth = 33.385416666666664
$ws.Columns(11).ColumnWidth = 23.944010416666668
$ws.Columns(12).ColumnWidth = 21.830729166666668

$ws.Range("D5").Select() | Out-Null
